# Daily attendance processing - 2025-12-29 15:02:07
# Swap the order of "Recorded By" names in column G:
#   "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"
# Only cells whose value is exactly the original combined string are touched;
# cells containing just "System" or just "dnasr281@gmail.com" are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column G is the "Recorded By" column.
$col = 7
$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2
    if ($val -eq $oldValue) {
        $cell.Value = $newValue
        $changed++
    }
}

Write-Host "Updated $changed cell(s) in column G."
